# draft-gandhi-ippm-stamp-srpm-02.pptx -- "Add files via upload" edit
#
# Applies the four textual changes captured in the commit's OOXML diff:
#   1. Handout master date field: 2/2/21 -> 2/3/21
#   2. Slide 3 ("Requirements and Scope"): merge the three runs that make up
#      the "High scale for number of test sessions and faster detection
#      interval" bullet into a single run (text itself is unchanged).
#   3. Slide 3: "RFC8972]" -> "RFC 8972]" (missing space before the number).
#   4. Slide 7 ("STAMP Return Path TLV - Usage"): "SR Path" -> "SR path"
#      (lower-case "path") in the "Bidir SR Path ..." bullet.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Handout master date placeholder field: 2/2/21 -> 2/3/21
# ---------------------------------------------------------------------
$hm = $p.HandoutMaster
$dateShape = $hm.Shapes.Item(2)
$dateShape.TextFrame.TextRange.Text = "2/3/21"

# ---------------------------------------------------------------------
# Slide 3: "Requirements and Scope" -- Content Placeholder 2
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange

# 2. "High scale for number " + "of test " + "sessions and faster
#    detection interval" (3 runs) -> single run with the same text.
#    Round-trip through a throwaway value first so the engine treats it
#    as a genuine content change and rebuilds the paragraph as one run
#    (setting the exact same text as already present is a no-op).
$highScalePara = $s3Body.Paragraphs(6)
$highScalePara.Text = "placeholder"
$highScalePara.Text = "High scale for number of test sessions and faster detection interval"

# 3. "RFC8972]" -> "RFC 8972]" inside the "STAMP Extensions [RFC8972]" line.
$scopePara = $s3Body.Paragraphs(10)
$scopeFull = $scopePara.Text
$rfcStart = $scopeFull.IndexOf("RFC8972]")
$rfcChars = $scopePara.Characters($rfcStart + 1, 8)
$rfcChars.Text = "RFC 8972]"

# ---------------------------------------------------------------------
# Slide 7: "STAMP Return Path TLV - Usage" -- Content Placeholder 2
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7Body = $s7.Shapes.Item(2).TextFrame.TextRange

# 4. " SR Path (forward and reverse) ..." -> " SR path (forward and
#    reverse) ..." -- only the second run (after "Bidir") changes.
$bidirPara = $s7Body.Paragraphs(2)
$bidirFull = $bidirPara.Text
$bidirRun1 = "Bidir"
$bidirRun2Old = " SR Path (forward and reverse) dynamically computed using CSPF by the head-end node"
$bidirRun2New = " SR path (forward and reverse) dynamically computed using CSPF by the head-end node"
$bidirChars = $bidirPara.Characters($bidirRun1.Length + 1, $bidirRun2Old.Length)
$bidirChars.Text = $bidirRun2New
